# Trade #17 closed at 2026-02-16 21:24:44 - leadlag DOWN +0.000%
# Append a new row (row 16) to the "leadlag" sheet with the new OPEN trade entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 16

$ws.Cells.Item($row, 1).Value = 17                 # Trade #
# Date column contains a date-like string ("2026-02-16") that must stay as TEXT
# (matching the rest of column B), so force text with a leading apostrophe to
# avoid Excel auto-converting it into a date serial number.
$ws.Cells.Item($row, 2).Value = "'2026-02-16"       # Date
$ws.Cells.Item($row, 3).Value = "21:24:44"          # Time
$ws.Cells.Item($row, 4).Value = "leadlag"           # Strategy
$ws.Cells.Item($row, 5).Value = "DOWN"              # Side
$ws.Cells.Item($row, 6).Value = 69273.17999999999   # Entry Price
$ws.Cells.Item($row, 7).Value = ""                  # Exit Price (blank, trade still open)
$ws.Cells.Item($row, 8).Value = "OPEN"              # Status
$ws.Cells.Item($row, 9).Value = 0                   # P&L %
$ws.Cells.Item($row, 10).Value = 0                  # P&L $
$ws.Cells.Item($row, 11).Value = 0.75               # Confidence
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.080% move"  # Entry Reason
$ws.Cells.Item($row, 13).Value = ""                 # Exit Reason (blank, trade still open)
$ws.Cells.Item($row, 14).Value = 0                  # Duration (min)
